$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H (Importe) amounts: remove thousands separator ".", convert decimal "," to "." ---
$ws.Range("H2:H282").NumberFormat = "@"
$ws.Range('H2').Value = '2400.00'
$ws.Range('H3').Value = '66265.00'
$ws.Range('H4').Value = '495000.00'
$ws.Range('H5').Value = '4965.00'
$ws.Range('H6').Value = '2134.00'
$ws.Range('H7').Value = '420.00'
$ws.Range('H8').Value = '3290.00'
$ws.Range('H9').Value = '21712.00'
$ws.Range('H10').Value = '1484.00'
$ws.Range('H11').Value = '4900.00'
$ws.Range('H12').Value = '19820.00'
$ws.Range('H13').Value = '100.00'
$ws.Range('H14').Value = '545600.00'
$ws.Range('H15').Value = '97884.00'
$ws.Range('H16').Value = '253198.32'
$ws.Range('H17').Value = '130532.50'
$ws.Range('H18').Value = '6597.27'
$ws.Range('H19').Value = '2640.00'
$ws.Range('H20').Value = '2392.00'
$ws.Range('H21').Value = '2839.00'
$ws.Range('H22').Value = '2700.00'
$ws.Range('H23').Value = '8654.50'
$ws.Range('H24').Value = '94325.75'
$ws.Range('H25').Value = '38513.00'
$ws.Range('H26').Value = '21909.44'
$ws.Range('H27').Value = '3327.50'
$ws.Range('H28').Value = '19550.00'
$ws.Range('H29').Value = '1725.00'
$ws.Range('H30').Value = '9875.00'
$ws.Range('H31').Value = '15000.00'
$ws.Range('H32').Value = '12100.00'
$ws.Range('H33').Value = '13160.00'
$ws.Range('H34').Value = '3000.00'
$ws.Range('H35').Value = '1050.00'
$ws.Range('H36').Value = '28292.00'
$ws.Range('H37').Value = '796.05'
$ws.Range('H38').Value = '144.48'
$ws.Range('H39').Value = '26811.63'
$ws.Range('H40').Value = '4164.45'
$ws.Range('H41').Value = '168.00'
$ws.Range('H42').Value = '8000.00'
$ws.Range('H43').Value = '3211.64'
$ws.Range('H44').Value = '1779.79'
$ws.Range('H45').Value = '13500.00'
$ws.Range('H46').Value = '26229.00'
$ws.Range('H47').Value = '5200.00'
$ws.Range('H48').Value = '535.41'
$ws.Range('H49').Value = '12322.00'
$ws.Range('H50').Value = '2350.00'
$ws.Range('H51').Value = '57765.03'
$ws.Range('H52').Value = '40.00'
$ws.Range('H53').Value = '14648.70'
$ws.Range('H54').Value = '63.03'
$ws.Range('H55').Value = '21051.30'
$ws.Range('H56').Value = '11800.00'
$ws.Range('H57').Value = '138462.69'
$ws.Range('H58').Value = '244.60'
$ws.Range('H59').Value = '890.00'
$ws.Range('H60').Value = '141031.98'
$ws.Range('H61').Value = '1128.20'
$ws.Range('H62').Value = '142936.78'
$ws.Range('H63').Value = '8500.00'
$ws.Range('H64').Value = '8855.00'
$ws.Range('H65').Value = '8418.15'
$ws.Range('H66').Value = '3200.00'
$ws.Range('H67').Value = '120.00'
$ws.Range('H68').Value = '860.00'
$ws.Range('H69').Value = '9313.15'
$ws.Range('H70').Value = '12946.97'
$ws.Range('H71').Value = '720.00'
$ws.Range('H72').Value = '7651.99'
$ws.Range('H73').Value = '8420.00'
$ws.Range('H74').Value = '62852.00'
$ws.Range('H75').Value = '44200.00'
$ws.Range('H76').Value = '4775.05'
$ws.Range('H77').Value = '65000.00'
$ws.Range('H78').Value = '6928.00'
$ws.Range('H79').Value = '11199.64'
$ws.Range('H80').Value = '1412.46'
$ws.Range('H81').Value = '7520.00'
$ws.Range('H82').Value = '950.00'
$ws.Range('H83').Value = '3600.00'
$ws.Range('H84').Value = '1480.00'
$ws.Range('H85').Value = '17750.00'
$ws.Range('H86').Value = '6600.00'
$ws.Range('H87').Value = '7350.00'
$ws.Range('H88').Value = '22255.00'
$ws.Range('H89').Value = '58410.00'
$ws.Range('H90').Value = '59200.00'
$ws.Range('H91').Value = '1315.00'
$ws.Range('H92').Value = '1591.00'
$ws.Range('H93').Value = '1700.00'
$ws.Range('H94').Value = '7000.00'
$ws.Range('H95').Value = '4777.00'
$ws.Range('H96').Value = '18834.41'
$ws.Range('H97').Value = '56.00'
$ws.Range('H98').Value = '1389.00'
$ws.Range('H99').Value = '1770.00'
$ws.Range('H100').Value = '5300.00'
$ws.Range('H101').Value = '10426.98'
$ws.Range('H102').Value = '3811.65'
$ws.Range('H103').Value = '1200.00'
$ws.Range('H104').Value = '508.00'
$ws.Range('H105').Value = '1190.00'
$ws.Range('H106').Value = '3840.00'
$ws.Range('H107').Value = '12000.00'
$ws.Range('H108').Value = '564.00'
$ws.Range('H109').Value = '120.00'
$ws.Range('H110').Value = '220.00'
$ws.Range('H111').Value = '11292.00'
$ws.Range('H112').Value = '20988.60'
$ws.Range('H113').Value = '26385.84'
$ws.Range('H114').Value = '300.00'
$ws.Range('H115').Value = '1401.60'
$ws.Range('H116').Value = '8950.00'
$ws.Range('H117').Value = '14985.00'
$ws.Range('H118').Value = '150.00'
$ws.Range('H119').Value = '2798.00'
$ws.Range('H120').Value = '4800.00'
$ws.Range('H121').Value = '9200.00'
$ws.Range('H122').Value = '1600.00'
$ws.Range('H123').Value = '1420.00'
$ws.Range('H124').Value = '8806.00'
$ws.Range('H125').Value = '9430.56'
$ws.Range('H126').Value = '2580.00'
$ws.Range('H127').Value = '25730.00'
$ws.Range('H128').Value = '17500.00'
$ws.Range('H129').Value = '3500.00'
$ws.Range('H130').Value = '6400.00'
$ws.Range('H131').Value = '3600.00'
$ws.Range('H132').Value = '20988.00'
$ws.Range('H133').Value = '66000.00'
$ws.Range('H134').Value = '32000.00'
$ws.Range('H135').Value = '2600.00'
$ws.Range('H136').Value = '6000.00'
$ws.Range('H137').Value = '137000.00'
$ws.Range('H138').Value = '2975.00'
$ws.Range('H139').Value = '1129.63'
$ws.Range('H140').Value = '2033.00'
$ws.Range('H141').Value = '7640.00'
$ws.Range('H142').Value = '11852.00'
$ws.Range('H143').Value = '10740.00'
$ws.Range('H144').Value = '3245.14'
$ws.Range('H145').Value = '10545.02'
$ws.Range('H146').Value = '50805.00'
$ws.Range('H147').Value = '14000.00'
$ws.Range('H148').Value = '25000.00'
$ws.Range('H149').Value = '16000.00'
$ws.Range('H150').Value = '7500.00'
$ws.Range('H151').Value = '36000.00'
$ws.Range('H152').Value = '8000.00'
$ws.Range('H153').Value = '13000.00'
$ws.Range('H154').Value = '10000.00'
$ws.Range('H155').Value = '17017.00'
$ws.Range('H156').Value = '7250.00'
$ws.Range('H157').Value = '12000.00'
$ws.Range('H158').Value = '14000.00'
$ws.Range('H159').Value = '11000.00'
$ws.Range('H160').Value = '12000.00'
$ws.Range('H161').Value = '12000.00'
$ws.Range('H162').Value = '10000.00'
$ws.Range('H163').Value = '8000.00'
$ws.Range('H164').Value = '12000.00'
$ws.Range('H165').Value = '14000.00'
$ws.Range('H166').Value = '7000.00'
$ws.Range('H167').Value = '36800.00'
$ws.Range('H168').Value = '12000.00'
$ws.Range('H169').Value = '16000.00'
$ws.Range('H170').Value = '13000.00'
$ws.Range('H171').Value = '1500.00'
$ws.Range('H172').Value = '6000.00'
$ws.Range('H173').Value = '16500.00'
$ws.Range('H174').Value = '20255.40'
$ws.Range('H175').Value = '3500.00'
$ws.Range('H176').Value = '30000.00'
$ws.Range('H177').Value = '13000.00'
$ws.Range('H178').Value = '12000.00'
$ws.Range('H179').Value = '75000.00'
$ws.Range('H180').Value = '7000.00'
$ws.Range('H181').Value = '50160.00'
$ws.Range('H182').Value = '18000.00'
$ws.Range('H183').Value = '29095.00'
$ws.Range('H184').Value = '14800.00'
$ws.Range('H185').Value = '200.00'
$ws.Range('H186').Value = '10500.00'
$ws.Range('H187').Value = '6700.00'
$ws.Range('H188').Value = '24500.00'
$ws.Range('H189').Value = '9700.00'
$ws.Range('H190').Value = '229.59'
$ws.Range('H191').Value = '5884.00'
$ws.Range('H192').Value = '14673.18'
$ws.Range('H193').Value = '2165.50'
$ws.Range('H194').Value = '39.00'
$ws.Range('H195').Value = '50541.00'
$ws.Range('H196').Value = '13330.00'
$ws.Range('H197').Value = '1174.96'
$ws.Range('H198').Value = '1183.37'
$ws.Range('H199').Value = '950.00'
$ws.Range('H200').Value = '1100.00'
$ws.Range('H201').Value = '1194.99'
$ws.Range('H202').Value = '1000.00'
$ws.Range('H203').Value = '27697.65'
$ws.Range('H204').Value = '1495.08'
$ws.Range('H205').Value = '11410.00'
$ws.Range('H206').Value = '16181.00'
$ws.Range('H207').Value = '10877.85'
$ws.Range('H208').Value = '530.00'
$ws.Range('H209').Value = '375.00'
$ws.Range('H210').Value = '6589.05'
$ws.Range('H211').Value = '6270.00'
$ws.Range('H212').Value = '170.62'
$ws.Range('H213').Value = '9178.00'
$ws.Range('H214').Value = '777.96'
$ws.Range('H215').Value = '4660.00'
$ws.Range('H216').Value = '57421.69'
$ws.Range('H217').Value = '2178.65'
$ws.Range('H218').Value = '6469.12'
$ws.Range('H219').Value = '15000.00'
$ws.Range('H220').Value = '30000.00'
$ws.Range('H221').Value = '30000.00'
$ws.Range('H222').Value = '63450.00'
$ws.Range('H223').Value = '30000.00'
$ws.Range('H224').Value = '30000.00'
$ws.Range('H225').Value = '30000.00'
$ws.Range('H226').Value = '60000.00'
$ws.Range('H227').Value = '60000.00'
$ws.Range('H228').Value = '60000.00'
$ws.Range('H229').Value = '30000.00'
$ws.Range('H230').Value = '39999.60'
$ws.Range('H231').Value = '35700.00'
$ws.Range('H232').Value = '4531034.27'
$ws.Range('H233').Value = '2700.00'
$ws.Range('H234').Value = '1800.00'
$ws.Range('H235').Value = '20000.00'
$ws.Range('H236').Value = '7294604.00'
$ws.Range('H237').Value = '165500.00'
$ws.Range('H238').Value = '159500.00'
$ws.Range('H239').Value = '150500.00'
$ws.Range('H240').Value = '163625.00'
$ws.Range('H241').Value = '150500.00'
$ws.Range('H242').Value = '150500.00'
$ws.Range('H243').Value = '269500.00'
$ws.Range('H244').Value = '337500.00'
$ws.Range('H245').Value = '394500.00'
$ws.Range('H246').Value = '150500.00'
$ws.Range('H247').Value = '199208.00'
$ws.Range('H248').Value = '150500.00'
$ws.Range('H249').Value = '150500.00'
$ws.Range('H250').Value = '150500.00'
$ws.Range('H251').Value = '269500.00'
$ws.Range('H252').Value = '388500.00'
$ws.Range('H253').Value = '269500.00'
$ws.Range('H254').Value = '150500.00'
$ws.Range('H255').Value = '278000.00'
$ws.Range('H256').Value = '150500.00'
$ws.Range('H257').Value = '150500.00'
$ws.Range('H258').Value = '155250.00'
$ws.Range('H259').Value = '150500.00'
$ws.Range('H260').Value = '471786.81'
$ws.Range('H261').Value = '8000.00'
$ws.Range('H262').Value = '85070.00'
$ws.Range('H263').Value = '453142.44'
$ws.Range('H264').Value = '500000.71'
$ws.Range('H265').Value = '26004.13'
$ws.Range('H266').Value = '289500.00'
$ws.Range('H267').Value = '110000.00'
$ws.Range('H268').Value = '70000.00'
$ws.Range('H269').Value = '1000.00'
$ws.Range('H270').Value = '1000.00'
$ws.Range('H271').Value = '20250.00'
$ws.Range('H272').Value = '20450.00'
$ws.Range('H273').Value = '9890.00'
$ws.Range('H274').Value = '55000.00'
$ws.Range('H275').Value = '6000.00'
$ws.Range('H276').Value = '49852.00'
$ws.Range('H277').Value = '25500.00'
$ws.Range('H278').Value = '22200.00'
$ws.Range('H279').Value = '1550.00'
$ws.Range('H280').Value = '2239.00'
$ws.Range('H281').Value = '600.00'
$ws.Range('H282').Value = '26481.00'
$ws.Range("H2:H282").Style = "Normal"

# --- Proveedor names: comma separators replaced with periods ---
$ws.Range('E172').Value = 'PARPAGNOLI. PEDRO RICARDO'
$ws.Range('F172').Value = 'PARPAGNOLI. PEDRO RICARDO'
$ws.Range('E191').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('E205').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
